$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 25; this pushes current rows 25-31 down to 26-32
$ws.Rows.Item(25).Insert()

# Copy style (number format) of the date cell from row 26 (previously row 25) into new row 25's D cell
$ws.Range("D26").Copy()
$ws.Range("D25").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Populate the new row 25 with data
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 44609
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100112030
$ws.Range("G25").Value = "Poroto granado"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 26000
$ws.Range("L25").Value = 28000
$ws.Range("M25").Value = 27000
$ws.Range("N25").Value = "$/saco 25 kilos"
$ws.Range("O25").Value = "Región Metropolitana"
$ws.Range("P25").Value = 1080
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"
